# 813_AFMC_PGI_5337_102_91.docx - "Added last minute updates"
#
# 1. Add a paragraph border (space-only, no line) around the first
#    paragraph and widen its left indent from 120 -> 225 twips.
# 2. Update the merge-field id text in that paragraph's first run from
#    **ID__AFFARS_pgi_5337_topic_6__ID** to
#    **ID__AFFARS_AFMC_PGI_5337_102_91__ID**, and drop the trailing
#    run that held a single literal space.

$d = $word.ActiveDocument
$p = $d.Paragraphs(1)

# --- paragraph border (w:pBdr w:top/left/bottom/right w:space="5", no line) ---
$p.Borders.DistanceFromTop = 5
$p.Borders.DistanceFromLeft = 5
$p.Borders.DistanceFromBottom = 5
$p.Borders.DistanceFromRight = 5

# --- left indent 120 -> 225 twips (twips = points * 20) ---
$p.Format.LeftIndent = 225 / 20

# --- drop the trailing single-space run at the end of the paragraph ---
$paraEnd = $p.Range.End
$trailingSpace = $d.Range($paraEnd - 2, $paraEnd - 1)
if ($trailingSpace.Text -eq " ") {
    $trailingSpace.Delete()
}

# --- update the merge-field id text ---
$d.Content.Find.Execute("**ID__AFFARS_pgi_5337_topic_6__ID**", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_AFMC_PGI_5337_102_91__ID**", 2)
